$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C10: change numeric value from 18 to 1.0
$ws.Range("C10").Value = 1.0

# B11: change text from "1" to "R40"
$ws.Range("B11").Value = "R40"
